$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 3.75
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.95

# Row 3 updates
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.85
